$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Cell($fromAddr, $toAddr) {
    $src = $ws.Range($fromAddr)
    $dst = $ws.Range($toAddr)
    $src.Copy($dst)
}

# ---------------------------------------------------------------------------
# New weekly block "Feb 22 - 3.0" at rows 59-69, mirroring the "Feb 19 - 3.0"
# block at rows 46-56 (13 rows lower). Copy cell-by-cell (not whole ranges)
# so that cells which are blank/unset in the source stay blank/unset in the
# destination instead of Excel materializing an empty <c> for every column.
# ---------------------------------------------------------------------------

# Row 59 <- Row 46 ("Total" banner row, fully styled A:J)
foreach ($col in @("A","B","C","D","E","F","G","H","I","J")) {
    $f = "${col}46"
    $t = "${col}59"
    Copy-Cell $f $t
}
$ws.Range("A59").Value = "Feb 22 - 3.0"

# Row 60 <- Row 47 (Total figures)
Copy-Cell "B47" "B60"
Copy-Cell "E47" "E60"
$ws.Range("E60").Value = 29367
Copy-Cell "J49" "J60"
$ws.Range("J60").Value = "common/WhirlyGlobeLib/"

# New J47 cell (added alongside this edit)
Copy-Cell "J49" "J47"
$ws.Range("J47").Value = "common/WhirlyGlobeLib/"

# Row 62 <- Row 49 (iOS WhirlyGlobeLib banner row, fully styled A:J)
foreach ($col in @("A","B","C","D","E","F","G","H","I","J")) {
    $f = "${col}49"
    $t = "${col}62"
    Copy-Cell $f $t
}

# Row 63 <- Row 50 (iOS WhirlyGlobeLib figures)
Copy-Cell "B50" "B63"
Copy-Cell "E50" "E63"
$ws.Range("E63").Value = 3732

# Row 65 <- Row 52 (Android Component banner row, fully styled A:J)
foreach ($col in @("A","B","C","D","E","F","G","H","I","J")) {
    $f = "${col}52"
    $t = "${col}65"
    Copy-Cell $f $t
}

# Row 66 <- Row 53 (Android Component figures)
Copy-Cell "B53" "B66"
Copy-Cell "E53" "E66"
$ws.Range("E66").Value = 26645

# Row 68 <- Row 55 (grand total)
Copy-Cell "B55" "B68"
Copy-Cell "E55" "E68"
$ws.Range("E68").Formula = "=E60+E63+E66"

# Row 69 <- Row 56 (core ratio) -- J56 moves down to J69
Copy-Cell "B56" "B69"
Copy-Cell "E56" "E69"
$ws.Range("E69").Formula = "=E60/E68"
Copy-Cell "J56" "J69"
$ws.Range("J56").ClearContents()

# ---------------------------------------------------------------------------
# View state: the saved workbook had scrolled back up and selected E14.
# (topLeftCell scroll position isn't tracked by this host; selection is.)
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
